$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add four new script filenames into column A for rows 3-6.
# (Row 3 already has data in B:E; rows 4 gets a new A value; rows 5-6 are brand new rows.)
$ws.Range("A3").Value = "SCRIPT/D73P11A/us0211.ssb"
$ws.Range("A4").Value = "SCRIPT/D73P11A/us0311.ssb"
$ws.Range("A5").Value = "SCRIPT/D73P11A/us0408.ssb"
$ws.Range("A6").Value = "SCRIPT/D73P11A/us2108.ssb"

# Give row 3 (B3:E3) a thin bottom border to visually separate it from row 4,
# since row 3 now marks the end of the first translated-line group.
$ws.Range("B3").Borders.Item(9).LineStyle = 1
$ws.Range("C3:E3").Borders.Item(9).LineStyle = 1

# Rows 4-6 now wrap their (8pt) text in column A onto two lines, so grow them
# to the same height already used for row 2's two-line wrapped entry.
$ws.Rows.Item(4).RowHeight = 43.2
$ws.Rows.Item(5).RowHeight = 43.2
$ws.Rows.Item(6).RowHeight = 43.2

# Leave the selection on the last originally-populated cell.
$ws.Range("E4").Select()
